$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new rows with data (columns F and G), continuing the existing sequence
$ws.Range("F9").Value = 7
$ws.Range("G9").Value = "projectile_create"

$ws.Range("F10").Value = 8
$ws.Range("G10").Value = "projectile_die"

# Update the selected cell to match the diff
$ws.Range("G10").Select()
